$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Revise previously reported weekly "waargenomen" (observed) deaths in column G ---
$ws.Range("G4").Value  = 3613
$ws.Range("G7").Value  = 4978
$ws.Range("G15").Value = 2682
$ws.Range("G21").Value = 2523
$ws.Range("G22").Value = 2669
$ws.Range("G23").Value = 2651
$ws.Range("G24").Value = 2628
$ws.Range("G25").Value = 3197
$ws.Range("G26").Value = 2821

# --- Insert a new row for week 35 right after week 34 (row 26), pushing the "Som" row down ---
$ws.Rows("27").Insert() | Out-Null

$ws.Range("F27").Value = 35
$ws.Range("G27").Value = 2689
$ws.Range("H27").Value = 2822
$ws.Range("I27").Formula = "=G27-H27"

# --- Extend the totals (now on row 29) to include the new week 35 row ---
$ws.Range("G29").Formula = "=SUM(G3:G27)"
$ws.Range("H29").Formula = "=SUM(H3:H27)"
$ws.Range("I29").Formula = "=SUM(I3:I27)"

# --- Match the saved view: scrolled down a bit, with the new total row selected ---
$ws.Application.ActiveWindow.ScrollRow = 17
$ws.Range("G29:I29").Select() | Out-Null
